$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.080397090533681
$ws.Range("D2").Value = 1.081885094789044
$ws.Range("E2").Value = 1.084132546534265
$ws.Range("F2").Value = 1.09085830725353
$ws.Range("I2").Value = 1.061277256252109
$ws.Range("J2").Value = 1.085276829674535
$ws.Range("K2").Value = 1.084554752430571
$ws.Range("L2").Value = 1.086796359566906
$ws.Range("M2").Value = 1.093504784812176
$ws.Range("N2").Value = 1.086818047255665
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.081901435684936
$ws.Range("D3").Value = 1.0831202368343
$ws.Range("E3").Value = 1.085589918653132
$ws.Range("F3").Value = 1.092226158267356
$ws.Range("I3").Value = 1.061821317675753
$ws.Range("J3").Value = 1.086439932652983
$ws.Range("K3").Value = 1.085607415161026
$ws.Range("L3").Value = 1.088071142368297
$ws.Range("M3").Value = 1.094691524689314
$ws.Range("N3").Value = 1.087982801973752
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.082873402459816
$ws.Range("D4").Value = 1.083918055467445
$ws.Range("E4").Value = 1.086531811769951
$ws.Range("F4").Value = 1.093110105613551
$ws.Range("I4").Value = 1.062171220273791
$ws.Range("J4").Value = 1.087190630320577
$ws.Range("K4").Value = 1.086286589017214
$ws.Range("L4").Value = 1.088894362888205
$ws.Range("M4").Value = 1.095457735714367
$ws.Range("N4").Value = 1.088734565718141
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.08328167856347
$ws.Range("D5").Value = 1.084253128271207
$ws.Range("E5").Value = 1.086927520924199
$ws.Range("F5").Value = 1.093481448970253
$ws.Range("I5").Value = 1.062317809486395
$ws.Range("J5").Value = 1.0875057724113
$ws.Range("K5").Value = 1.086571647157218
$ws.Range("L5").Value = 1.089240055841553
$ws.Range("M5").Value = 1.095779451761374
$ws.Range("N5").Value = 1.08905015534678
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.083350210264284
$ws.Range("D6").Value = 1.08430936926037
$ws.Range("E6").Value = 1.086993947012195
$ws.Range("F6").Value = 1.093543783597034
$ws.Range("I6").Value = 1.062342392645373
$ws.Range("J6").Value = 1.087558659855412
$ws.Range("K6").Value = 1.086619482456765
$ws.Range("L6").Value = 1.089298076571082
$ws.Range("M6").Value = 1.095833446052935
$ws.Range("N6").Value = 1.089103117897129
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.082878859187091
$ws.Range("D7").Value = 1.083922534017443
$ws.Range("E7").Value = 1.08653710027938
$ws.Range("F7").Value = 1.093115068568067
$ws.Range("I7").Value = 1.062173181005826
$ws.Range("J7").Value = 1.087194843030741
$ws.Range("K7").Value = 1.086290399805173
$ws.Range("L7").Value = 1.088898983575418
$ws.Range("M7").Value = 1.095462036061445
$ws.Range("N7").Value = 1.088738784410836
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.080905794611791
$ws.Range("D8").Value = 1.082302810033884
$ws.Range("E8").Value = 1.084625308877497
$ws.Range("F8").Value = 1.091320819164982
$ws.Range("I8").Value = 1.061461569491328
$ws.Range("J8").Value = 1.085670304103522
$ws.Range("K8").Value = 1.084910915452747
$ws.Range("L8").Value = 1.087227523400259
$ws.Range("M8").Value = 1.093906202026862
$ws.Range("N8").Value = 1.087212080463474
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.077417609657249
$ws.Range("D9").Value = 1.079437685772858
$ws.Range("E9").Value = 1.081247591410162
$ws.Range("F9").Value = 1.08815010653601
$ws.Range("I9").Value = 1.060191089577746
$ws.Range("J9").Value = 1.082969013248545
$ws.Range("K9").Value = 1.082464784354582
$ws.Range("L9").Value = 1.084269302454855
$ws.Range("M9").Value = 1.091151449904067
$ws.Range("N9").Value = 1.08450695346564
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.075084036557334
$ws.Range("D10").Value = 1.077519901165445
$ws.Range("E10").Value = 1.078989393989854
$ws.Range("F10").Value = 1.086029875614793
$ws.Range("I10").Value = 1.059332809152863
$ws.Range("J10").Value = 1.081157814635211
$ws.Range("K10").Value = 1.080823430467259
$ws.Range("L10").Value = 1.082288114867107
$ws.Range("M10").Value = 1.089305756579179
$ws.Range("N10").Value = 1.082693182742527
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.07407154642895
$ws.Range("D11").Value = 1.076687578307957
$ws.Range("E11").Value = 1.07800996441811
$ws.Range("F11").Value = 1.085110187955056
$ws.Range("I11").Value = 1.058958444967629
$ws.Range("J11").Value = 1.080371015558477
$ws.Range("K11").Value = 1.080110122249814
$ws.Range("L11").Value = 1.081428011438875
$ws.Range("M11").Value = 1.088504295706209
$ws.Range("N11").Value = 1.081905266320858
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.073695147516551
$ws.Range("D12").Value = 1.076378123182923
$ws.Range("E12").Value = 1.077645910119893
$ws.Range("F12").Value = 1.084768325512289
$ws.Range("I12").Value = 1.058818976763353
$ws.Range("J12").Value = 1.080078375214934
$ws.Range("K12").Value = 1.079844772431785
$ws.Range("L12").Value = 1.081108187531941
$ws.Range("M12").Value = 1.088206250832635
$ws.Range("N12").Value = 1.081612210394461
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.073775900758436
$ws.Range("D13").Value = 1.076444515737291
$ws.Range("E13").Value = 1.077724012467584
$ws.Range("F13").Value = 1.084841667570908
$ws.Range("I13").Value = 1.058848911898598
$ws.Range("J13").Value = 1.080141165235765
$ws.Range("K13").Value = 1.079901708858885
$ws.Range("L13").Value = 1.081176806492935
$ws.Range("M13").Value = 1.088270198260853
$ws.Range("N13").Value = 1.081675089584323
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.074040439654842
$ws.Range("D14").Value = 1.076662004672831
$ws.Range("E14").Value = 1.077979876733114
$ws.Range("F14").Value = 1.085081934629571
$ws.Range("I14").Value = 1.058946924921546
$ws.Range("J14").Value = 1.080346833773774
$ws.Range("K14").Value = 1.080088196452246
$ws.Range("L14").Value = 1.081401581738174
$ws.Range("M14").Value = 1.088479666333033
$ws.Range("N14").Value = 1.081881050195246
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.074203388793024
$ws.Range("D15").Value = 1.076795967716766
$ws.Range("E15").Value = 1.078137489715123
$ws.Range("F15").Value = 1.085229937821398
$ws.Range("I15").Value = 1.059007259189061
$ws.Range("J15").Value = 1.080473501311863
$ws.Range("K15").Value = 1.080203045037818
$ws.Range("L15").Value = 1.081540027490668
$ws.Range("M15").Value = 1.088608680398082
$ws.Range("N15").Value = 1.082007897615768
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.075151188507106
$ws.Range("D16").Value = 1.077575098909974
$ws.Range("E16").Value = 1.079054360803198
$ws.Range("F16").Value = 1.086090877656081
$ws.Range("I16").Value = 1.059357596794851
$ws.Range("J16").Value = 1.08120997784318
$ws.Range("K16").Value = 1.080870715272627
$ws.Range("L16").Value = 1.082345149369125
$ws.Range("M16").Value = 1.089358898650154
$ws.Range("N16").Value = 1.082745420028235
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.07574516646363
$ws.Range("D17").Value = 1.078063311544075
$ws.Range("E17").Value = 1.07962905183428
$ws.Range("F17").Value = 1.086630485472644
$ws.Range("I17").Value = 1.059576622761221
$ws.Range("J17").Value = 1.081671265793235
$ws.Range("K17").Value = 1.081288829157064
$ws.Range("L17").Value = 1.082849577278712
$ws.Range("M17").Value = 1.089828880125562
$ws.Range("N17").Value = 1.083207363060089
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.076091427762935
$ws.Range("D18").Value = 1.078347893918704
$ws.Range("E18").Value = 1.079964104083607
$ws.Range("F18").Value = 1.086945074366186
$ws.Range("I18").Value = 1.059704114404715
$ws.Range("J18").Value = 1.081940082820943
$ws.Range("K18").Value = 1.081532458317798
$ws.Range("L18").Value = 1.083143586282756
$ws.Range("M18").Value = 1.09010279442142
$ws.Range("N18").Value = 1.083476561838811
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.076209460941203
$ws.Range("D19").Value = 1.078444898174408
$ws.Range("E19").Value = 1.080078322136813
$ws.Range("F19").Value = 1.087052314934265
$ws.Range("I19").Value = 1.059747541343637
$ws.Range("J19").Value = 1.082031701226114
$ws.Range("K19").Value = 1.081615487413619
$ws.Range("L19").Value = 1.083243799461302
$ws.Range("M19").Value = 1.090196155347464
$ws.Range("N19").Value = 1.083568310352625
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.075681458587114
$ws.Range("D20").Value = 1.078010950006632
$ws.Range("E20").Value = 1.079567409046683
$ws.Range("F20").Value = 1.086572606766711
$ws.Range("I20").Value = 1.059553150525028
$ws.Range("J20").Value = 1.08162179929363
$ws.Range("K20").Value = 1.08124399535516
$ws.Range("L20").Value = 1.082795479255913
$ws.Range("M20").Value = 1.089778478151997
$ws.Range("N20").Value = 1.083157826312381
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.073962548310367
$ws.Range("D21").Value = 1.076597967752451
$ws.Range("E21").Value = 1.077904538064054
$ws.Range("F21").Value = 1.085011188913585
$ws.Range("I21").Value = 1.058918073949852
$ws.Range("J21").Value = 1.080286280284821
$ws.Range("K21").Value = 1.08003329146784
$ws.Range("L21").Value = 1.081335400563201
$ws.Range("M21").Value = 1.088417992827951
$ws.Range("N21").Value = 1.081820410713394
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.072879974276459
$ws.Range("D22").Value = 1.075707869212168
$ws.Range("E22").Value = 1.076857571549815
$ws.Range("F22").Value = 1.084028016490507
$ws.Range("I22").Value = 1.05851638669409
$ws.Range("J22").Value = 1.079444337444002
$ws.Range("K22").Value = 1.079269783248498
$ws.Range("L22").Value = 1.080415401877246
$ws.Range("M22").Value = 1.087560592986139
$ws.Range("N22").Value = 1.080977272217168
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.073454043403871
$ws.Range("D23").Value = 1.076179890654432
$ws.Range("E23").Value = 1.07741272873158
$ws.Range("F23").Value = 1.084549354327212
$ws.Range("I23").Value = 1.058729556392375
$ws.Range("J23").Value = 1.079890882741924
$ws.Range("K23").Value = 1.079674752572832
$ws.Range("L23").Value = 1.080903301578102
$ws.Range("M23").Value = 1.088015309469454
$ws.Range("N23").Value = 1.081424451660633
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.075710246040096
$ws.Range("D24").Value = 1.078034610499718
$ws.Range("E24").Value = 1.079595263249186
$ws.Range("F24").Value = 1.086598760142468
$ws.Range("I24").Value = 1.059563757429447
$ws.Range("J24").Value = 1.081644151831044
$ws.Range("K24").Value = 1.081264254591849
$ws.Range("L24").Value = 1.08281992448919
$ws.Range("M24").Value = 1.089801253308023
$ws.Range("N24").Value = 1.083180210592962
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.078320786026684
$ws.Range("D25").Value = 1.080179722194668
$ws.Range("E25").Value = 1.082121908088314
$ws.Range("F25").Value = 1.088970917168651
$ws.Range("I25").Value = 1.060521516406877
$ws.Range("J25").Value = 1.083669158177244
$ws.Range("K25").Value = 1.083099013262059
$ws.Range("L25").Value = 1.085035639562294
$ws.Range("M25").Value = 1.091865214825236
$ws.Range("N25").Value = 1.085208092680445
